$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "212.84" into numbers),
# then restore the default "Normal" style so no stray formatting remains.
function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '26.250.42'
Set-TextValue $ws.Range("E2") '  +0.55%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.607.64'
Set-TextValue $ws.Range("E3") '  +0.63%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.01%  '

# Row 5
Set-TextValue $ws.Range("D5") '212.84'
Set-TextValue $ws.Range("E5") '  +0.38%  '

# Row 6
Set-TextValue $ws.Range("E6") '  -0.08%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +0.33%  '

# Row 8
Set-TextValue $ws.Range("E8") '  +0.44%  '

# Row 10
Set-TextValue $ws.Range("D10") '18.21'
Set-TextValue $ws.Range("E10") '  +1.78%  '

# Row 11
Set-TextValue $ws.Range("E11") '  -0.81%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.832.46'

# Row 13
Set-TextValue $ws.Range("D13") '1.620.14'
Set-TextValue $ws.Range("E13") '  +1.36%  '

# Row 14
Set-TextValue $ws.Range("E14") '  +0.61%  '

# Row 15
Set-TextValue $ws.Range("E15") '  +1.14%  '

# Row 16
Set-TextValue $ws.Range("D16") '26.298.51'
Set-TextValue $ws.Range("E16") '  +0.81%  '

# Row 17
Set-TextValue $ws.Range("D17") '62.12'
Set-TextValue $ws.Range("E17") '  +2.85%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.0₃0728'
Set-TextValue $ws.Range("E18") '  +1.11%  '

# Row 19
Set-TextValue $ws.Range("E19") '  -0.08%  '

# Row 20
Set-TextValue $ws.Range("D20") '202.59'
Set-TextValue $ws.Range("E20") '  -0.74%  '

# Row 21
Set-TextValue $ws.Range("E21") '  +1.15%  '

# Row 22
Set-TextValue $ws.Range("E22") '  +0.16%  '

# Row 23
Set-TextValue $ws.Range("D23") '6.01'
Set-TextValue $ws.Range("E23") '  +0.84%  '

# Row 24
Set-TextValue $ws.Range("E24") '  +2.99%  '

# Row 25
Set-TextValue $ws.Range("D25") '144.92'
Set-TextValue $ws.Range("E25") '  +2.15%  '

# Row 26
Set-TextValue $ws.Range("E26") '  -0.05%  '

# Row 27
Set-TextValue $ws.Range("E27") '  -3.43%  '

# Row 28
Set-TextValue $ws.Range("E28") '  +0.18%  '

# Row 29
Set-TextValue $ws.Range("E29") '  +2.06%  '

# Row 30
Set-TextValue $ws.Range("D30") '0.0493'
Set-TextValue $ws.Range("E30") '  +5.12%  '

# Row 31
Set-TextValue $ws.Range("E31") '  +0.77%  '

# Row 32
Set-TextValue $ws.Range("E32") '  +2.92%  '

# Row 33
Set-TextValue $ws.Range("E33") '  -2.03%  '

# Row 34
Set-TextValue $ws.Range("E34") '  +2.71%  '

# Row 35
Set-TextValue $ws.Range("E35") '  +1.19%  '

# Row 36
Set-TextValue $ws.Range("D36") '1.165.83'
Set-TextValue $ws.Range("E36") '  +5.46%  '

# Row 37
Set-TextValue $ws.Range("E37") '  +1.51%  '

# Row 38
Set-TextValue $ws.Range("E38") '  -0.14%  '

# Row 39
Set-TextValue $ws.Range("E39") '  +0.17%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.788'
Set-TextValue $ws.Range("E40") '  +1.72%  '

# Row 41
Set-TextValue $ws.Range("E41") '  +1.05%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.783'
Set-TextValue $ws.Range("E42") '  +0.83%  '

# Row 43
Set-TextValue $ws.Range("E43") '  +3.70%  '

# Row 44
Set-TextValue $ws.Range("D44") '1.745.04'
Set-TextValue $ws.Range("E44") '  +0.50%  '

# Row 45
Set-TextValue $ws.Range("D45") '92.05'
Set-TextValue $ws.Range("E45") '  -0.58%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +0.93%  '

# Row 47
Set-TextValue $ws.Range("D47") '54.21'
Set-TextValue $ws.Range("E47") '  +1.70%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.0₇0996'
Set-TextValue $ws.Range("E48") '  -3.46%  '

# Row 49
Set-TextValue $ws.Range("E49") '  +0.34%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.408'
Set-TextValue $ws.Range("E50") '  -0.37%  '

# Row 51
Set-TextValue $ws.Range("E51") '  -0.11%  '
